$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$v = $ws.Range("C7").Value
Write-Host "value is: $v"
Write-Host ("value is: " + $ws.Range("C7").Value.ToString())
